# Automatische test-sync: 2025-06-19 21:36:50
# Adds the new "Probleem met inloggen" mail-log entry (row 19) to the
# "Logs" sheet, extends the conditional-formatting ranges to cover it,
# and refreshes the "Dashboard" category summary so the new IT entry is
# reflected (count 1 -> 2) and the rows are kept sorted by count desc.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new incoming mail as row 19
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A19").Value = "Probleem met inloggen"
$logs.Range("B19").Value = "mailmind.test@zohomail.eu"
$logs.Range("C19").Value = "Ik kan niet inloggen op mijn account. Kunnen jullie dit oplossen?"
$logs.Range("D19").Value = "IT / Technisch probleem"
$logs.Range("F19").Value = "2025-06-19 21:36:10"
$logs.Range("G19").Value = "Nee"

# Extend the conditional-formatting ranges (D2:D18 -> D2:D19, G2:G18 -> G2:G19)
# so the newly added row is covered by the same highlighting rules.
$catFormats = $logs.Range("D2:D18").FormatConditions
$catFormats.Item(1).ModifyAppliesToRange($logs.Range("D2:D19"))

$answeredFormats = $logs.Range("G2:G18").FormatConditions
$answeredFormats.Item(1).ModifyAppliesToRange($logs.Range("G2:G19"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: refresh the category counts / ordering
#    (sorted by "Aantal" descending; the new mail bumps
#    "IT / Technisch probleem" from 1 to 2, moving it above
#    "Offerte / Prijsaanvraag")
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "IT / Technisch probleem"
$dash.Range("B4").Value = 2
$dash.Range("A5").Value = "Offerte / Prijsaanvraag"
$dash.Range("B5").Value = 2
$dash.Range("A6").Value = "Factuur / Administratie"
$dash.Range("B6").Value = 1
$dash.Range("A7").Value = "Openingstijden / Locatie"
$dash.Range("B7").Value = 1
